$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new value would otherwise be
# auto-detected by Excel as a number, so they remain plain text like the original.
$textCells = @('D5', 'D6', 'D9', 'D11', 'D12', 'D13', 'D14', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D29', 'D30', 'D31', 'D32', 'D35', 'D36', 'D37', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '62.206.73'
$ws.Range('E2').Value = '  -2.80%  '
$ws.Range('D3').Value = '3.177.60'
$ws.Range('E3').Value = '  -4.34%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '584.95'
$ws.Range('E5').Value = '  -2.45%  '
$ws.Range('D6').Value = '134.73'
$ws.Range('E6').Value = '  -6.42%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.175.83'
$ws.Range('E8').Value = '  -4.26%  '
$ws.Range('D9').Value = '0.502'
$ws.Range('E9').Value = '  -4.38%  '
$ws.Range('E10').Value = '  -6.20%  '
$ws.Range('D11').Value = '5.24'
$ws.Range('E11').Value = '  -5.63%  '
$ws.Range('D12').Value = '0.449'
$ws.Range('E12').Value = '  -5.54%  '
$ws.Range('D13').Value = '0.0000234'
$ws.Range('E13').Value = '  -6.84%  '
$ws.Range('D14').Value = '33.13'
$ws.Range('E14').Value = '  -5.05%  '
$ws.Range('D15').Value = '3.702.50'
$ws.Range('E15').Value = '  -4.26%  '
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('D17').Value = '3.173.84'
$ws.Range('E17').Value = '  -4.44%  '
$ws.Range('D18').Value = '62.280.54'
$ws.Range('E18').Value = '  -2.81%  '
$ws.Range('D19').Value = '6.56'
$ws.Range('E19').Value = '  -5.30%  '
$ws.Range('D20').Value = '454.14'
$ws.Range('E20').Value = '  -6.22%  '
$ws.Range('D21').Value = '13.91'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').Value = '0.702'
$ws.Range('E22').Value = '  -5.26%  '
$ws.Range('D23').Value = '7.58'
$ws.Range('E23').Value = '  -5.66%  '
$ws.Range('D24').Value = '13.32'
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').Value = '82.20'
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D29').Value = '6.88'
$ws.Range('E29').Value = '  -5.61%  '
$ws.Range('D30').Value = '7.79'
$ws.Range('E30').Value = '  -5.73%  '
$ws.Range('D31').Value = '2.01'
$ws.Range('E31').Value = '  -7.11%  '
$ws.Range('D32').Value = '27.18'
$ws.Range('E32').Value = '  -8.79%  '
$ws.Range('E33').Value = '  -4.21%  '
$ws.Range('E34').Value = '  -7.31%  '
$ws.Range('D35').Value = '1.04'
$ws.Range('E35').Value = '  -6.29%  '
$ws.Range('D36').Value = '5.77'
$ws.Range('E36').Value = '  -4.28%  '
$ws.Range('D37').Value = '51.01'
$ws.Range('E37').Value = '  -4.46%  '
$ws.Range('D38').Value = '0.0₃0685'
$ws.Range('E38').Value = '  -10.16%  '
$ws.Range('E39').Value = '  -4.56%  '
$ws.Range('D40').Value = '2.943.52'
$ws.Range('E40').Value = '  -3.68%  '
$ws.Range('D41').Value = '407.43'
$ws.Range('E41').Value = '  -6.20%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Value = '8.00'
$ws.Range('E42').Value = '  -5.44%  '
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '2.62'
$ws.Range('E44').Value = '  -6.35%  '
$ws.Range('D45').Value = '0.249'
$ws.Range('E45').Value = '  -7.26%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').Value = '2.13'
$ws.Range('E46').Value = '  -4.25%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').Value = '0.999'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').Value = '35.72'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('D49').Value = '25.36'
$ws.Range('E49').Value = '  -4.83%  '
$ws.Range('D50').Value = '123.07'
$ws.Range('E50').Value = '  -0.58%  '
